$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 12, shifting rows 12:39
# down to 13:40 (this also grows the used range from A1:R39 to A1:R40).
$ws.Rows.Item(12).EntireRow.Insert()

# Populate the newly-inserted row 12 with the new market record.
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "Vega Monumental Concepción"
$ws.Range("C12").Value = "Bíobío"
$ws.Range("D12").Value = 44645
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = 100112030
$ws.Range("G12").Value = "Poroto granado"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 220
$ws.Range("K12").Value = 26000
$ws.Range("L12").Value = 27000
$ws.Range("M12").Value = 26455
$ws.Range("N12").Value = "$/saco 25 kilos"
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 1058
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"
